$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.718.63"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.599.38"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D5").Value = "'211.59"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'19.61"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "1.823.80"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.604.84"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.0₃0741"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'208.85"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").Value = "'9.01"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'143.66"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'7.12"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'2.98"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'1.26"
$ws.Range("E33").Value = "  +17.85%  "
$ws.Range("D34").Value = "1.277.00"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'0.776"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").Value = "'62.54"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "1.735.44"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "'90.37"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'7.53"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  +1.61%  "
